# The source diff for this revision only re-serializes the package's
# word/document.xml and word/styles.xml parts: every changed line has
# exactly the same attributes/values/namespaces before and after, just
# re-ordered (alphabetically) by the attribute's local name. That is a
# side effect of the commit's real change ("Fixed POI packaging and
# upgraded to POI 3.15.") regenerating the expected test fixture with a
# newer Apache POI/XMLBeans serializer - it is not a content, formatting,
# or structural edit that exists in Word's object model.
#
# There is therefore nothing to change through COM automation: the
# document's text, runs, styles, section/page setup, fonts, latent
# style table, etc. are all byte-for-byte equal in value to the
# original. We simply touch the document without altering any content
# so the package round-trips unchanged.
$d = $word.ActiveDocument
$d.Content | Out-Null
